# Insert a new weekly price record at row 316 (Terminal Hortofrutícola Agro
# Chillán - Pepino ensalada, Hortaliza). Existing rows 316:346 shift down to
# 317:347.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing row 316 (and everything below it) down by one row.
$ws.Rows.Item(316).Insert()

# Populate the newly-inserted row 316 with the new observation.
$ws.Range("A316").Value = 7
$ws.Range("B316").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C316").Value = "Ñuble"
$ws.Range("D316").Value = 45106
$ws.Range("E316").Value = 16
$ws.Range("F316").Value = 100112043
$ws.Range("G316").Value = "Pepino ensalada"
$ws.Range("H316").Value = "Sin especificar"
$ws.Range("I316").Value = "Primera"
$ws.Range("J316").Value = 110
$ws.Range("K316").Value = 14000
$ws.Range("L316").Value = 15000
$ws.Range("M316").Value = 14455
$ws.Range("N316").Value = "$/caja 60 unidades"
$ws.Range("O316").Value = "Región de Arica y Parinacota"
$ws.Range("P316").Value = 241
$ws.Range("Q316").Value = 60
$ws.Range("R316").Value = "Hortaliza"
